$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header labels: "<Field>_old" -> "<Field>_FV2410" and
#    "<Field>_new" -> "<Field>_FV2504". Column K ("diff") stays as-is.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Segmentname_FV2410"
$ws.Range("B1").Value = "Segmentgruppe_FV2410"
$ws.Range("C1").Value = "Segment_FV2410"
$ws.Range("D1").Value = "Datenelement_FV2410"
$ws.Range("E1").Value = "Segment ID_FV2410"
$ws.Range("F1").Value = "Code_FV2410"
$ws.Range("G1").Value = "Qualifier_FV2410"
$ws.Range("H1").Value = "Beschreibung_FV2410"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("J1").Value = "Bedingung_FV2410"

$ws.Range("L1").Value = "Segmentname_FV2504"
$ws.Range("M1").Value = "Segmentgruppe_FV2504"
$ws.Range("N1").Value = "Segment_FV2504"
$ws.Range("O1").Value = "Datenelement_FV2504"
$ws.Range("P1").Value = "Segment ID_FV2504"
$ws.Range("Q1").Value = "Code_FV2504"
$ws.Range("R1").Value = "Qualifier_FV2504"
$ws.Range("S1").Value = "Beschreibung_FV2504"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2504"
$ws.Range("U1").Value = "Bedingung_FV2504"

# ---------------------------------------------------------------------------
# 2) Freeze the header row (pane split under row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3) Turn the used range into an Excel Table (ListObject).
#    Excel's default "Insert Table" always bolds the header row and records
#    that as a new dxf/tableStyleInfo-style entry. The header row here
#    already carries its own explicit formatting (fill/border/bold) that
#    must stay untouched, so: clear the header formatting right before
#    adding the table (so Excel sees nothing to "diff" into a new dxf),
#    then Undo() just that formatting change (Undo only rolls back the
#    single last action, the table itself is unaffected) to restore the
#    original header look without leaving a dxf behind.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U80")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

$excel.Undo()

$tbl.TableStyle = ""
